# Add 5 new daily data rows (334-338) to the COVID tracking sheet, and
# extend the used range accordingly (A1:R333 -> A1:R338).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates formatted as plain text (e.g. "2021/02/25").
# Force text formatting first so Excel doesn't auto-convert the
# strings into date serial numbers.
$ws.Range("A334:A338").NumberFormat = "@"

function Set-DataRow {
    param($r, $vals)

    $ws.Range("A$r").Value = $vals[0]

    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
    $ws.Range("E$r").Value = $vals[4]
    $ws.Range("F$r").Value = $vals[5]
    $ws.Range("G$r").Value = $vals[6]

    if ($vals[7] -eq $null) {
        $ws.Range("H$r").Value = ""
    } else {
        $ws.Range("H$r").Value = $vals[7]
    }

    $ws.Range("I$r").Value = $vals[8]
    $ws.Range("J$r").Value = $vals[9]
    $ws.Range("K$r").Value = $vals[10]
    $ws.Range("L$r").Value = $vals[11]
    $ws.Range("M$r").Value = $vals[12]
    $ws.Range("N$r").Value = $vals[13]
    $ws.Range("O$r").Value = $vals[14]
    $ws.Range("P$r").Value = $vals[15]
    $ws.Range("Q$r").Value = $vals[16]
    $ws.Range("R$r").Value = $vals[17]
}

# r, A..R
Set-DataRow 334 @("2021/02/25", 333, 8222, 172, 135, 0.02091948431038677, 7892, 18556, 26778, 0,   419, 3, 416, 0,      35, 25, 0, 48)
Set-DataRow 335 @("2021/02/26", 334, 8222, 172, 135, 0.02091948431038677, 7892, 18556, 26778, 0,   419, 3, 416, 0,      35, 25, 0, 48)
Set-DataRow 336 @("2021/02/27", 335, 8222, 172, 135, 0.02091948431038677, 7892, 18556, 26778, 0,   419, 3, 416, 0,      35, 25, 0, 48)
Set-DataRow 337 @("2021/02/28", 336, 8222, 172, 135, 0.02091948431038677, 7892, 18556, 26778, 0,   419, 3, 416, 0,      35, 25, 0, 49)
Set-DataRow 338 @("2021/03/01", 337, 8447, 177, 83,  0.02095418491772227, 8163, $null, 8447,  225, 502, 1, 501, -18331, 32, 31, 5, 49)

# Drop back to the default "Normal" style so no extra number formatting
# is left behind on the new cells (matches the rest of column A).
$ws.Range("A334:A338").Style = "Normal"
